$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend the table with a new "2021" column (column R), matching the
# formatting of the corresponding cells from the rest of the table.
$ws.Range("Q4").Copy()
$ws.Range("R4").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("R4").Value = 2021

$ws.Range("D5").Copy()
$ws.Range("R5").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("R5").Value = 31.8

$ws.Range("Q6").Copy()
$ws.Range("R6").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("R6").Value = 12957.1

$excel.CutCopyMode = 0

# Update the active cell selection to match the new edit area
$ws.Range("R4:R6").Select()
